# Adapt column header formatting to respective input file names (#7)
#
# - Rename the two header-row suffixes used throughout the sheet:
#     "..._old" -> "..._FV2310"
#     "..._new" -> "..._FV2404"
# - Turn the A1:U88 range into a native Excel Table ("Table1") that carries
#   the (renamed) header row as its column names, with autofilter.
# - Freeze the header row (row 1) so it stays visible while scrolling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header cells in row 1 (A1:U1) -------------------------------
# Every header ends in either "_old" or "_new" (except the plain "diff"
# column) - swap those generic suffixes for the concrete format versions
# being compared.

$lastCol = 21
for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $current = $cell.Value()
    if ($null -ne $current) {
        $updated = $current -replace "_old$", "_FV2310"
        $updated = $updated -replace "_new$", "_FV2404"
        if ($updated -ne $current) {
            $cell.Value = $updated
        }
    }
}

# --- 2. Wrap the used range in a native table ------------------------------

$tableRange = $ws.Range("A1:U88")
$table = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$table.Name = "Table1"

# --- 3. Freeze the header row ----------------------------------------------

[void]$ws.Activate()
[void]$ws.Range("A2").Select()
[void]($excel.ActiveWindow.FreezePanes = $true)
